# Weekly update: add two new daily price rows (Primera/Segunda) for
# "Terminal Hortofrutícola Agro Chillán - Brócoli" dated 2023-04-25 (serial 45041).
# This pushes all existing data rows from row 442 onward down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows right before the current row 442, shifting the
# rest of the table (rows 442..462) down to rows 444..464.
$ws.Rows.Item(442).Insert()
$ws.Rows.Item(442).Insert()

# --- New row 442: "Primera" quality record for the new date ---
$ws.Cells.Item(442, 1).Value = 7
$ws.Cells.Item(442, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(442, 3).Value = "Ñuble"
$ws.Cells.Item(442, 4).Value = 45041
$ws.Cells.Item(442, 5).Value = 16
$ws.Cells.Item(442, 6).Value = 100112023
$ws.Cells.Item(442, 7).Value = "Brócoli"
$ws.Cells.Item(442, 8).Value = "Sin especificar"
$ws.Cells.Item(442, 9).Value = "Primera"
$ws.Cells.Item(442, 10).Value = 500
$ws.Cells.Item(442, 11).Value = 1200
$ws.Cells.Item(442, 12).Value = 1200
$ws.Cells.Item(442, 13).Value = 1200
$ws.Cells.Item(442, 14).Value = "$/unidad"
$ws.Cells.Item(442, 15).Value = "Región del Maule"
$ws.Cells.Item(442, 16).Value = 1200
$ws.Cells.Item(442, 17).Value = 1
$ws.Cells.Item(442, 18).Value = "Hortaliza"

# --- New row 443: "Segunda" quality record for the same new date ---
$ws.Cells.Item(443, 1).Value = 7
$ws.Cells.Item(443, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(443, 3).Value = "Ñuble"
$ws.Cells.Item(443, 4).Value = 45041
$ws.Cells.Item(443, 5).Value = 16
$ws.Cells.Item(443, 6).Value = 100112023
$ws.Cells.Item(443, 7).Value = "Brócoli"
$ws.Cells.Item(443, 8).Value = "Sin especificar"
$ws.Cells.Item(443, 9).Value = "Segunda"
$ws.Cells.Item(443, 10).Value = 300
$ws.Cells.Item(443, 11).Value = 1000
$ws.Cells.Item(443, 12).Value = 1000
$ws.Cells.Item(443, 13).Value = 1000
$ws.Cells.Item(443, 14).Value = "$/unidad"
$ws.Cells.Item(443, 15).Value = "Región del Maule"
$ws.Cells.Item(443, 16).Value = 1000
$ws.Cells.Item(443, 17).Value = 1
$ws.Cells.Item(443, 18).Value = "Hortaliza"

# Make sure the date column keeps using the existing date number format
# (same style the rest of column D already uses).
$ws.Range("D442").NumberFormat = $ws.Range("D444").NumberFormat()
$ws.Range("D443").NumberFormat = $ws.Range("D444").NumberFormat()
